$d = $word.ActiveDocument

# Locate the paragraph built from several runs split around spell-check
# proofErr markers ("Selaamla" / " be " / "yenı" / " update").
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Selaamla*") {
        $target = $p
    }
}

$r = $target.Range

# Replace the whole paragraph -- including its own end-of-paragraph mark, so
# the orphaned w:proofErr elements go away with it -- with clean OOXML:
#   1) the same wording merged into a single plain run (no proofErr markers)
#   2) a new blank paragraph
#   3) a new paragraph with the extra sentence, carrying the _GoBack bookmark
#      that used to sit at the end of the original paragraph
$xmlFragment = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Selaamla be yenı update</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Selamlar asdghsdlgkhsldghdslgh şlsdhg</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

[void]$r.InsertXML($xmlFragment)
